$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-like string into a cell as literal text (not an
# auto-converted date serial), while keeping the cell's existing style
# index ("s=2") intact. We briefly force a text number format so Excel's
# smart input parsing leaves the value alone, then restore the original
# "General" formatting by copying it over from a neighboring cell that
# keeps style 2 (D4, which stays "Maths"/General) - this only touches the
# display format, not the already-committed text value.
function Set-TextValue($range, [string]$text) {
  $range.NumberFormat = "@"
  $range.Value = $text
  $ws.Range("D4").Copy() | Out-Null
  $range.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# Row 4 updates
Set-TextValue $ws.Range("A4") "06-30-2020"
$ws.Range("B4").Value = "Mr. B P Tripathi"
$ws.Range("C4").Value = "III B"
$ws.Range("E4").Value = 95
$ws.Range("F4").Value = 56
$ws.Range("G4").Value = 39
$ws.Range("H4").Value = "Give and Take"
$ws.Range("J4").Value = "Quiz"
$ws.Range("K4").Value = ""

# Row 5 updates
Set-TextValue $ws.Range("A5") "06-30-2020"
$ws.Range("B5").Value = "Mr. Ashok Uttam"
$ws.Range("C5").Value = "V B"
$ws.Range("D5").Value = "CS"
$ws.Range("E5").Value = 95
$ws.Range("F5").Value = 32
$ws.Range("G5").Value = 63
$ws.Range("H5").Value = "OOPs"
$ws.Range("J5").Value = "Sth"
$ws.Range("K5").Value = "Keep It Up"

$excel.CutCopyMode = $false

# Remove rows 6-8 entirely (report now only covers two entries)
$ws.Rows("6:8").Delete()
